$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 0

$ws.Range("H8").Select() | Out-Null
